# Add a new worksheet "severalAnswers3" at the end of the workbook.
# It is a duplicate of "severalAnswers" (same identifiers data/layout as
# "severalAnswers2", but carrying the same cell styling already baked
# into "severalAnswers"/"severalAnswers1").
$wb = $excel.ActiveWorkbook

# Remember which sheet was active so we can restore the selection after
# the copy (Excel activates the newly inserted sheet by default).
$originalActiveSheet = $wb.ActiveSheet

$sourceSheet = $wb.Worksheets.Item("severalAnswers")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# "Move or Copy" -> create a copy, placed after the current last sheet.
$sourceSheet.Copy([Type]::Missing, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "severalAnswers3"

$originalActiveSheet.Activate()
